$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 540.1667
$ws.Cells.Item(4, 9).Value = 348.2
$ws.Cells.Item(4, 11).Value = 348.2
$ws.Cells.Item(4, 13).Value = -234.2

# Row 12
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 13).ClearContents()

# Row 107
$ws.Cells.Item(107, 8).Value = 132.11765
$ws.Cells.Item(107, 9).Value = 143.71428
$ws.Cells.Item(107, 11).Value = 143.71428
$ws.Cells.Item(107, 13).Value = 1776.28572

# Row 110
$ws.Cells.Item(110, 8).Value = 47701
$ws.Cells.Item(110, 10).Value = 47701
$ws.Cells.Item(110, 12).Value = 47701
$ws.Cells.Item(110, 14).Value = -55881

# Row 112
$ws.Cells.Item(112, 8).Value = 1189.129
$ws.Cells.Item(112, 9).Value = 869.25
$ws.Cells.Item(112, 11).Value = 2607.75
$ws.Cells.Item(112, 13).Value = -1499.75

# Row 120
$ws.Cells.Item(120, 8).Value = 46873
$ws.Cells.Item(120, 10).Value = 46873
$ws.Cells.Item(120, 12).Value = 46873
$ws.Cells.Item(120, 14).Value = -56549

# Row 135
$ws.Cells.Item(135, 8).Value = 1581.5862
$ws.Cells.Item(135, 9).Value = 1437.2307
$ws.Cells.Item(135, 11).Value = 12935.0763
$ws.Cells.Item(135, 13).Value = -10400.0763

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 655.1875
$ws.Cells.Item(2, 9).Value = 575.75
$ws.Cells.Item(2, 11).Value = 575.75
$ws.Cells.Item(2, 13).Value = -462.75

# Row 52
$ws.Cells.Item(52, 8).Value = 56986.5
$ws.Cells.Item(52, 10).Value = 56986.5
$ws.Cells.Item(52, 12).Value = 56986.5
$ws.Cells.Item(52, 14).Value = -57622.5

# Row 104
$ws.Cells.Item(104, 8).Value = 18063.875
$ws.Cells.Item(104, 10).Value = 18063.875
$ws.Cells.Item(104, 12).Value = 18063.875
$ws.Cells.Item(104, 14).Value = -25051.875

# Row 115
$ws.Cells.Item(115, 8).Value = 39998.77
$ws.Cells.Item(115, 10).Value = 39998.77
$ws.Cells.Item(115, 12).Value = 39998.77
$ws.Cells.Item(115, 14).Value = -43132.77

# Row 116
$ws.Cells.Item(116, 8).Value = 655.1875
$ws.Cells.Item(116, 9).Value = 575.75
$ws.Cells.Item(116, 11).Value = 575.75
$ws.Cells.Item(116, 13).Value = 1718.25

# Row 135
$ws.Cells.Item(135, 8).Value = 64398.8
$ws.Cells.Item(135, 10).Value = 64398.8
$ws.Cells.Item(135, 12).Value = 64398.8
$ws.Cells.Item(135, 14).Value = -74538.8

# Row 139
$ws.Cells.Item(139, 8).Value = 110353
$ws.Cells.Item(139, 10).Value = 110353
$ws.Cells.Item(139, 12).Value = 110353
$ws.Cells.Item(139, 14).Value = -120633

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 655.1875
$ws.Cells.Item(3, 9).Value = 575.75
$ws.Cells.Item(3, 11).Value = 575.75
$ws.Cells.Item(3, 13).Value = -461.75

# Row 22
$ws.Cells.Item(22, 8).Value = 14654292
$ws.Cells.Item(22, 9).Value = 14654292
$ws.Cells.Item(22, 11).Value = 14654292
$ws.Cells.Item(22, 13).Value = -14654119

# Row 52
$ws.Cells.Item(52, 8).Value = 99989.2
$ws.Cells.Item(52, 10).Value = 99989.2
$ws.Cells.Item(52, 12).Value = 99989.2
$ws.Cells.Item(52, 14).Value = -100515.2

# Row 53
$ws.Cells.Item(53, 8).Value = 38997.332
$ws.Cells.Item(53, 10).Value = 38997.332
$ws.Cells.Item(53, 12).Value = 38997.332
$ws.Cells.Item(53, 14).Value = -40145.332

# Row 81
$ws.Cells.Item(81, 8).Value = 25270.334
$ws.Cells.Item(81, 10).Value = 25270.334
$ws.Cells.Item(81, 12).Value = 25270.334
$ws.Cells.Item(81, 14).Value = -27392.334

# Row 84
$ws.Cells.Item(84, 8).Value = 25270.334
$ws.Cells.Item(84, 10).Value = 25270.334
$ws.Cells.Item(84, 12).Value = 75811.00199999999
$ws.Cells.Item(84, 14).Value = -86419.00199999999

# Row 109
$ws.Cells.Item(109, 8).Value = 66853.42999999999
$ws.Cells.Item(109, 10).Value = 66853.42999999999
$ws.Cells.Item(109, 12).Value = 66853.42999999999
$ws.Cells.Item(109, 14).Value = -69627.42999999999

# Row 118
$ws.Cells.Item(118, 8).Value = 74781
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 13).ClearContents()

# Row 119
$ws.Cells.Item(119, 8).Value = 99986
$ws.Cells.Item(119, 10).Value = 99986
$ws.Cells.Item(119, 12).Value = 99986
$ws.Cells.Item(119, 14).Value = -109662

# Row 121
$ws.Cells.Item(121, 8).Value = 99989.2
$ws.Cells.Item(121, 10).Value = 99989.2
$ws.Cells.Item(121, 12).Value = 99989.2
$ws.Cells.Item(121, 14).Value = -103483.2

# Row 122
$ws.Cells.Item(122, 8).Value = 77712.60000000001
$ws.Cells.Item(122, 10).Value = 77712.60000000001
$ws.Cells.Item(122, 12).Value = 77712.60000000001
$ws.Cells.Item(122, 14).Value = -87512.60000000001

# Row 132
$ws.Cells.Item(132, 8).Value = 47045.812
$ws.Cells.Item(132, 10).Value = 47045.812
$ws.Cells.Item(132, 12).Value = 47045.812
$ws.Cells.Item(132, 14).Value = -57165.812

# Row 135
$ws.Cells.Item(135, 8).Value = 119499.5
$ws.Cells.Item(135, 10).Value = 119499.5
$ws.Cells.Item(135, 12).Value = 119499.5
$ws.Cells.Item(135, 14).Value = -129639.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1862.1389
$ws.Cells.Item(31, 9).Value = 1506.9706
$ws.Cells.Item(31, 10).Value = 7900
$ws.Cells.Item(31, 11).Value = 1506.9706
$ws.Cells.Item(31, 12).Value = 7900
$ws.Cells.Item(31, 13).Value = -1211.9706
$ws.Cells.Item(31, 14).Value = -8490

# Row 34
$ws.Cells.Item(34, 8).Value = 1862.1389
$ws.Cells.Item(34, 9).Value = 1506.9706
$ws.Cells.Item(34, 10).Value = 7900
$ws.Cells.Item(34, 11).Value = 1506.9706
$ws.Cells.Item(34, 12).Value = 7900
$ws.Cells.Item(34, 13).Value = -1304.9706
$ws.Cells.Item(34, 14).Value = -8304

# Row 86
$ws.Cells.Item(86, 8).Value = 8657.5
$ws.Cells.Item(86, 10).Value = 8657.5
$ws.Cells.Item(86, 12).Value = 8657.5
$ws.Cells.Item(86, 14).Value = -10903.5

# Row 89
$ws.Cells.Item(89, 8).Value = 8657.5
$ws.Cells.Item(89, 10).Value = 8657.5
$ws.Cells.Item(89, 12).Value = 43287.5
$ws.Cells.Item(89, 14).Value = -54519.5

# Row 138
$ws.Cells.Item(138, 8).Value = 105234.664
$ws.Cells.Item(138, 10).Value = 124997.5
$ws.Cells.Item(138, 12).Value = 124997.5
$ws.Cells.Item(138, 14).Value = -135277.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 67.92308
$ws.Cells.Item(2, 9).Value = 109
$ws.Cells.Item(2, 10).Value = 20
$ws.Cells.Item(2, 11).Value = 654
$ws.Cells.Item(2, 12).Value = 120
$ws.Cells.Item(2, 13).Value = -541
$ws.Cells.Item(2, 14).Value = -346

# Row 54
$ws.Cells.Item(54, 8).Value = 3923.75
$ws.Cells.Item(54, 9).Value = 700
$ws.Cells.Item(54, 10).Value = 4998.3335
$ws.Cells.Item(54, 11).Value = 2100
$ws.Cells.Item(54, 12).Value = 14995.0005
$ws.Cells.Item(54, 13).Value = -1541
$ws.Cells.Item(54, 14).Value = -16113.0005

# Row 88
$ws.Cells.Item(88, 8).Value = 4068.1667
$ws.Cells.Item(88, 10).Value = 4681.8
$ws.Cells.Item(88, 12).Value = 14045.4
$ws.Cells.Item(88, 14).Value = -14901.4

# Row 91
$ws.Cells.Item(91, 8).Value = 4068.1667
$ws.Cells.Item(91, 10).Value = 4681.8
$ws.Cells.Item(91, 12).Value = 14045.4
$ws.Cells.Item(91, 14).Value = -17009.4

# Row 138
$ws.Cells.Item(138, 8).Value = 6603.281
$ws.Cells.Item(138, 9).Value = 6830.067
$ws.Cells.Item(138, 10).Value = 5752.8335
$ws.Cells.Item(138, 11).Value = 20490.201
$ws.Cells.Item(138, 12).Value = 17258.5005
$ws.Cells.Item(138, 13).Value = -15350.201
$ws.Cells.Item(138, 14).Value = -27538.5005

# Row 141
$ws.Cells.Item(141, 8).Value = 1629.4546
$ws.Cells.Item(141, 9).Value = 1629.4546
$ws.Cells.Item(141, 11).Value = 4888.3638
$ws.Cells.Item(141, 13).Value = 291.6361999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 4064.8
$ws.Cells.Item(22, 9).Value = 4479.6
$ws.Cells.Item(22, 10).Value = 3650
$ws.Cells.Item(22, 11).Value = 4479.6
$ws.Cells.Item(22, 12).Value = 3650
$ws.Cells.Item(22, 13).Value = -3950.6
$ws.Cells.Item(22, 14).Value = -4708

# Row 102
$ws.Cells.Item(102, 8).Value = 2251
$ws.Cells.Item(102, 9).Value = 2119.6365
$ws.Cells.Item(102, 11).Value = 2119.6365
$ws.Cells.Item(102, 13).Value = -497.6365000000001

# Row 107
$ws.Cells.Item(107, 8).Value = 1008.619
$ws.Cells.Item(107, 10).Value = 1029
$ws.Cells.Item(107, 12).Value = 1029
$ws.Cells.Item(107, 14).Value = -4869

# Row 113
$ws.Cells.Item(113, 8).Value = 3033.5
$ws.Cells.Item(113, 9).Value = 2733.3333
$ws.Cells.Item(113, 11).Value = 2733.3333
$ws.Cells.Item(113, 13).Value = -563.3332999999998

# Row 114
$ws.Cells.Item(114, 8).Value = 72648.55
$ws.Cells.Item(114, 10).Value = 72648.55
$ws.Cells.Item(114, 12).Value = 72648.55
$ws.Cells.Item(114, 14).Value = -81326.55

# Row 135
$ws.Cells.Item(135, 8).Value = 53519.523
$ws.Cells.Item(135, 10).Value = 53519.523
$ws.Cells.Item(135, 12).Value = 53519.523
$ws.Cells.Item(135, 14).Value = -63659.523

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1800
$ws.Cells.Item(16, 9).Value = 1800
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1800
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 14).Value = -1630

# Row 61
$ws.Cells.Item(61, 8).Value = 1844.9231
$ws.Cells.Item(61, 9).Value = 1780.3636
$ws.Cells.Item(61, 11).Value = 1780.3636
$ws.Cells.Item(61, 13).Value = -1578.3636

# Row 93
$ws.Cells.Item(93, 8).Value = 1335.8334
$ws.Cells.Item(93, 9).Value = 1335.8334
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 1335.8334
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).ClearContents()
$ws.Cells.Item(93, 14).Value = -87.83339999999998

# Row 94
$ws.Cells.Item(94, 8).Value = 30000
$ws.Cells.Item(94, 10).Value = 30000
$ws.Cells.Item(94, 12).Value = 30000
$ws.Cells.Item(94, 14).Value = -31352

# Row 113
$ws.Cells.Item(113, 8).Value = 1844.9231
$ws.Cells.Item(113, 9).Value = 1780.3636
$ws.Cells.Item(113, 11).Value = 1780.3636
$ws.Cells.Item(113, 13).Value = 389.6364000000001

# Row 118
$ws.Cells.Item(118, 8).Value = 50872.727
$ws.Cells.Item(118, 10).Value = 51960
$ws.Cells.Item(118, 12).Value = 51960
$ws.Cells.Item(118, 14).Value = -55274

$ws = $wb.Worksheets.Item("WVR")
# Row 121
$ws.Cells.Item(121, 8).Value = 40348.4
$ws.Cells.Item(121, 10).Value = 40348.4
$ws.Cells.Item(121, 12).Value = 40348.4
$ws.Cells.Item(121, 14).Value = -43842.4

Write-Host "Applied all updates"